# Weekly price-list update: a new daily record is inserted at the top of the
# data block (row 18), pushing all the subsequent records down by one row.
#
# Before: rows 18-31 held the historical records.
# After:  a brand-new record (Fecha 44741 / Perfection / ... ) occupies row 18,
#         and the old rows 18-31 are shifted down to become rows 19-32.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 18; this shifts rows 18:31 down to 19:32,
# automatically carrying their existing values/formatting with them.
$ws.Rows("18:18").Insert()

# Populate the newly inserted row 18 with the new weekly record.
$ws.Range("A18").Value = 11
$ws.Range("B18").Value = 'Vega Monumental Concepción'
$ws.Range("C18").Value = 'Bíobío'
$ws.Range("D18").Value = 44741
$ws.Range("E18").Value = 8
$ws.Range("F18").Value = 100112022
$ws.Range("G18").Value = 'Arveja Verde'
$ws.Range("H18").Value = 'Perfection'
$ws.Range("I18").Value = 'Primera'
$ws.Range("J18").Value = 100
$ws.Range("K18").Value = 40000
$ws.Range("L18").Value = 42000
$ws.Range("M18").Value = 41000
$ws.Range("N18").Value = '$/saco 25 kilos'
$ws.Range("O18").Value = 'Provincia de Limarí'
$ws.Range("P18").Value = 1640
$ws.Range("Q18").Value = 25
$ws.Range("R18").Value = 'Hortaliza'
